$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values on rows 5 and 6
$ws.Range("B5").Value = 215
$ws.Range("E5").Value = 4
$ws.Range("B6").Value = 225
$ws.Range("E6").Value = 5

# Move the active selection to E6 to match the saved view state
$ws.Range("E6").Select()
